$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40"; change it to the text "1".
# Use NumberFormat "@" (text) to ensure the new value is stored as a
# shared string (t="s") rather than being interpreted as a numeric value.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
